# Auto-generated Excel COM-interop edit script
# Applies updated market-price data cells across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 54.75
$ws.Range("I6").Value = 54.75
$ws.Range("K6").Value = 164.25
$ws.Range("M6").Value = -52.25
$ws.Range("H9").Value = 273.625
$ws.Range("I9").Value = 63
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 63
$ws.Range("L9").Value = 400
$ws.Range("M9").Value = 106
$ws.Range("N9").Value = -738
$ws.Range("H12").Value = 1731.7142
$ws.Range("I12").Value = 1853.5
$ws.Range("J12").Value = 1001
$ws.Range("K12").Value = 1853.5
$ws.Range("L12").Value = 1001
$ws.Range("M12").Value = -1683.5
$ws.Range("N12").Value = -1341
$ws.Range("H29").Value = 3342.8572
$ws.Range("J29").Value = 3733.3333
$ws.Range("L29").Value = 11199.9999
$ws.Range("N29").Value = -11761.9999
$ws.Range("H42").Value = 50
$ws.Range("I42").Value = 50
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 150
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 80
$ws.Range("N42").ClearContents()
$ws.Range("H55").Value = 515.0833
$ws.Range("I55").Value = 753
$ws.Range("J55").Value = 182
$ws.Range("K55").Value = 753
$ws.Range("L55").Value = 182
$ws.Range("M55").Value = -539
$ws.Range("N55").Value = -610
$ws.Range("H86").Value = 9324.846
$ws.Range("I86").Value = 2150.5
$ws.Range("J86").Value = 15474.286
$ws.Range("K86").Value = 2150.5
$ws.Range("L86").Value = 15474.286
$ws.Range("M86").Value = -1027.5
$ws.Range("N86").Value = -17720.286
$ws.Range("H89").Value = 9324.846
$ws.Range("I89").Value = 2150.5
$ws.Range("J89").Value = 15474.286
$ws.Range("K89").Value = 10752.5
$ws.Range("L89").Value = 77371.42999999999
$ws.Range("M89").Value = -5136.5
$ws.Range("N89").Value = -88603.42999999999
$ws.Range("H100").Value = 2714.4285
$ws.Range("I100").Value = 2300.4
$ws.Range("K100").Value = 2300.4
$ws.Range("M100").Value = -1759.4
$ws.Range("H112").Value = 1039.2941
$ws.Range("J112").Value = 1085.0646
$ws.Range("L112").Value = 3255.1938
$ws.Range("N112").Value = -5471.1938
$ws.Range("H113").Value = 66670076
$ws.Range("I113").Value = 100002420
$ws.Range("J113").Value = 5400
$ws.Range("K113").Value = 100002420
$ws.Range("L113").Value = 5400
$ws.Range("M113").Value = -99999166
$ws.Range("N113").Value = -11908
$ws.Range("H132").Value = 3078.5
$ws.Range("I132").Value = 3196.4348
$ws.Range("J132").Value = 2536
$ws.Range("K132").Value = 9589.304400000001
$ws.Range("L132").Value = 7608
$ws.Range("M132").Value = -7059.304400000001
$ws.Range("N132").Value = -12668
$ws.Range("H137").Value = 73517.5
$ws.Range("I137").Value = 2190.6667
$ws.Range("J137").Value = 127012.625
$ws.Range("K137").Value = 6572.000100000001
$ws.Range("L137").Value = 381037.875
$ws.Range("M137").Value = -4022.000100000001
$ws.Range("N137").Value = -386137.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3061.149
$ws.Range("I45").Value = 2402.1428
$ws.Range("K45").Value = 2402.1428
$ws.Range("M45").Value = -2025.1428
$ws.Range("H110").Value = 348.6842
$ws.Range("I110").Value = 302.64706
$ws.Range("J110").Value = 740
$ws.Range("K110").Value = 302.64706
$ws.Range("L110").Value = 740
$ws.Range("M110").Value = 1742.35294
$ws.Range("N110").Value = -4830
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1518.5769
$ws.Range("I99").Value = 1155.1666
$ws.Range("J99").Value = 2336.25
$ws.Range("K99").Value = 1155.1666
$ws.Range("L99").Value = 2336.25
$ws.Range("M99").Value = 342.8334
$ws.Range("N99").Value = -5332.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1599.3334
$ws.Range("I16").Value = 1905.5
$ws.Range("K16").Value = 1905.5
$ws.Range("M16").Value = -1618.5
$ws.Range("H31").Value = 2333.24
$ws.Range("I31").Value = 1249.4117
$ws.Range("J31").Value = 4636.375
$ws.Range("K31").Value = 1249.4117
$ws.Range("L31").Value = 4636.375
$ws.Range("M31").Value = -954.4117000000001
$ws.Range("N31").Value = -5226.375
$ws.Range("H34").Value = 2333.24
$ws.Range("I34").Value = 1249.4117
$ws.Range("J34").Value = 4636.375
$ws.Range("K34").Value = 1249.4117
$ws.Range("L34").Value = 4636.375
$ws.Range("M34").Value = -1047.4117
$ws.Range("N34").Value = -5040.375
$ws.Range("H99").Value = 29171656
$ws.Range("I99").Value = 7579730.5
$ws.Range("J99").Value = 55561788
$ws.Range("K99").Value = 7579730.5
$ws.Range("L99").Value = 55561788
$ws.Range("M99").Value = -7578232.5
$ws.Range("N99").Value = -55564784
$ws.Range("H113").Value = 1599.3334
$ws.Range("I113").Value = 1905.5
$ws.Range("K113").Value = 1905.5
$ws.Range("M113").Value = 264.5
$ws.Range("H122").Value = 1207.2858
$ws.Range("I122").Value = 1215.5385
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 3646.6155
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -1196.6155
$ws.Range("N122").Value = -8200
$ws.Range("H126").Value = 29171656
$ws.Range("I126").Value = 7579730.5
$ws.Range("J126").Value = 55561788
$ws.Range("K126").Value = 22739191.5
$ws.Range("L126").Value = 166685364
$ws.Range("M126").Value = -22736721.5
$ws.Range("N126").Value = -166690304
$ws.Range("H134").Value = 1028.6129
$ws.Range("I134").Value = 817.2632
$ws.Range("J134").Value = 1363.25
$ws.Range("K134").Value = 2451.7896
$ws.Range("L134").Value = 4089.75
$ws.Range("M134").Value = 83.21039999999994
$ws.Range("N134").Value = -9159.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 9000
$ws.Range("I105").Value = 9000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 27000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -24379
$ws.Range("N105").ClearContents()
$ws.Range("H131").Value = 773.6799999999999
$ws.Range("J131").Value = 788.70215
$ws.Range("L131").Value = 2366.10645
$ws.Range("N131").Value = -12446.10645

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4818986.5
$ws.Range("I70").Value = 22362.6
$ws.Range("J70").Value = 7816876
$ws.Range("K70").Value = 22362.6
$ws.Range("L70").Value = 7816876
$ws.Range("M70").Value = -22092.6
$ws.Range("N70").Value = -7817416
$ws.Range("H73").Value = 4818986.5
$ws.Range("I73").Value = 22362.6
$ws.Range("J73").Value = 7816876
$ws.Range("K73").Value = 22362.6
$ws.Range("L73").Value = 7816876
$ws.Range("M73").Value = -21426.6
$ws.Range("N73").Value = -7818748
$ws.Range("H107").Value = 109.53846
$ws.Range("I107").Value = 123.22222
$ws.Range("K107").Value = 123.22222
$ws.Range("M107").Value = 1796.77778
$ws.Range("H113").Value = 3410
$ws.Range("I113").Value = 2516.6667
$ws.Range("J113").Value = 4750
$ws.Range("K113").Value = 2516.6667
$ws.Range("L113").Value = 4750
$ws.Range("M113").Value = -346.6667000000002
$ws.Range("N113").Value = -9090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 817
$ws.Range("I46").Value = 501
$ws.Range("J46").Value = 975
$ws.Range("K46").Value = 501
$ws.Range("L46").Value = 975
$ws.Range("N46").Value = -1351
$ws.Range("M46").Value = -313

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 8735.286
$ws.Range("I52").Value = 8633.333000000001
$ws.Range("J52").Value = 8811.75
$ws.Range("K52").Value = 8633.333000000001
$ws.Range("L52").Value = 8811.75
$ws.Range("M52").Value = -8407.333000000001
$ws.Range("N52").Value = -9263.75
$ws.Range("H100").Value = 800
$ws.Range("I100").Value = 620.2
$ws.Range("J100").Value = 1249.5
$ws.Range("K100").Value = 1240.4
$ws.Range("L100").Value = 2499
$ws.Range("M100").Value = -699.4000000000001
$ws.Range("N100").Value = -3581
$ws.Range("H113").Value = 1352102.5
$ws.Range("I113").Value = 848.13336
$ws.Range("K113").Value = 2544.40008
$ws.Range("M113").Value = -374.4000800000003
